$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric remain stored as text,
# matching the source data which represents prices as formatted strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.489.81'
$ws.Range("E2").Value = '  -1.88%  '
$ws.Range("D3").Value = '1.748.37'
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '328.43'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '0.4782'
$ws.Range("E7").Value = '  +5.95%  '
$ws.Range("D8").Value = '0.3529'
$ws.Range("E8").Value = '  -1.98%  '
$ws.Range("D9").Value = '42.72'
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").Value = '0.07447'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").Value = '1.073'
$ws.Range("E11").Value = '  -3.09%  '
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '20.30'
$ws.Range("E13").Value = '  -3.50%  '
$ws.Range("D14").Value = '6.024'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.755.69'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '7.048'
$ws.Range("E16").Value = '  -2.86%  '
$ws.Range("D17").Value = '92.30'
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").Value = '0.00001066'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '0.06406'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '16.66'
$ws.Range("E21").Value = '  -3.23%  '
$ws.Range("D22").Value = '5.750'
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("D23").Value = '27.553.54'
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").Value = '11.01'
$ws.Range("E24").Value = '  -3.65%  '
$ws.Range("D25").Value = '2.149'
$ws.Range("E25").Value = '  +3.07%  '
$ws.Range("D26").Value = '161.97'
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("D27").Value = '19.93'
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").Value = '1.958.38'
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("D29").Value = '2.185'
$ws.Range("E29").Value = '  -2.54%  '
$ws.Range("D30").Value = '121.49'
$ws.Range("E30").Value = '  -3.83%  '
$ws.Range("D31").Value = '1.049'
$ws.Range("E31").Value = '  -5.38%  '
$ws.Range("D32").Value = '0.09371'
$ws.Range("E32").Value = '  +1.70%  '
$ws.Range("D33").Value = '3.648'
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").Value = '5.474'
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("D35").Value = '0.02251'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").Value = '11.31'
$ws.Range("E36").Value = '  -5.86%  '
$ws.Range("D37").Value = '0.05929'
$ws.Range("E37").Value = '  -3.40%  '
$ws.Range("D38").Value = '0.2046'
$ws.Range("E38").Value = '  -2.81%  '
$ws.Range("D39").Value = '4.830'
$ws.Range("E39").Value = '  -3.27%  '
$ws.Range("B40").Value = 'WEMIXTOKEN'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '1.433'
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.174'
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6073'
$ws.Range("E42").Value = '  -4.77%  '
$ws.Range("D43").Value = '7.720'
$ws.Range("E43").Value = '  -3.62%  '
$ws.Range("B44").Value = 'PancakeSwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D44").Value = '3.728'
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '12.79'
$ws.Range("E45").Value = '  -4.58%  '
$ws.Range("D46").Value = '0.5706'
$ws.Range("E46").Value = '  -4.00%  '
$ws.Range("D47").Value = '122.40'
$ws.Range("E47").Value = '  -0.79%  '
$ws.Range("D48").Value = '1.905'
$ws.Range("E48").Value = '  -3.26%  '
$ws.Range("D49").Value = '1.131'
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("D50").Value = '0.06751'
$ws.Range("E50").Value = '  -2.87%  '
$ws.Range("D51").Value = '71.29'
$ws.Range("E51").Value = '  -2.41%  '
